$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 542, shifting existing rows 542-653 down to 543-654
$ws.Rows.Item(542).Insert()

# Populate the newly inserted row 542 with its data
$ws.Cells.Item(542, 1).Value = 6
$ws.Cells.Item(542, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(542, 3).Value = "Metropolitana"
$ws.Cells.Item(542, 4).Value = 45015
$ws.Cells.Item(542, 5).Value = 13
$ws.Cells.Item(542, 6).Value = 100112039
$ws.Cells.Item(542, 7).Value = "Ciboulette"
$ws.Cells.Item(542, 8).Value = "Sin especificar"
$ws.Cells.Item(542, 9).Value = "Primera"
$ws.Cells.Item(542, 10).Value = 670
$ws.Cells.Item(542, 11).Value = 900
$ws.Cells.Item(542, 12).Value = 1000
$ws.Cells.Item(542, 13).Value = 942
$ws.Cells.Item(542, 14).Value = "$/docena de atados"
$ws.Cells.Item(542, 15).Value = "Región Metropolitana"
$ws.Cells.Item(542, 16).Value = 314
$ws.Cells.Item(542, 17).Value = 3
$ws.Cells.Item(542, 18).Value = "Hortaliza"
